$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arbeitspakete")

# F7: add TODO note about checking client names
$ws.Range("F7").Value = "TODO: prüfe Client-Namen"

# C11: 0.8 -> 0.9
$ws.Range("C11").Value = 0.9

# Row 13: clientseitige Firewall-Configuration
$ws.Range("C13").Value = 0.9
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = "parameterisierbare Funktion zum Erstellen und Löschen von Firewall-Regeln exisitiert; inkl.UnitTest; LB-Kandidaten dürfen max ""Hauptbenutzer"" sein; Stolperstein: globaler Firewall-Service Zustand gelöst"

# Row 14: USB-Ports deaktivieren
$ws.Range("C14").Value = 0.9

# Row 28: Modale Progressbars?
$ws.Range("C28").Value = 0.9
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = "sauberes Threading mit QThreads und Workers (QRunnables in ThreadPools) implementiert, ProgressBar am unteren Bildschirmrand funktioniert ordentlich"
$ws.Rows.Item(28).RowHeight = 35.05

# Update view: scroll position and active cell selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F7").Select()
